$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like cells (Coin name / Link) - plain text, no special number format needed
$textCells = @{
    "B12" = "TigerCash"
    "C12" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B13" = "LEO"
    "C13" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B14" = "GateToken"
    "C14" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "B15" = "BTSEToken"
    "C15" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B16" = "One"
    "C16" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
}

# Numeric-looking text cells (Price / Volume) - must force Text format to preserve exact string & avoid numeric coercion
$numericTextCells = @{
    "D2" = "246.54"
    "E2" = "0.58%"
    "D3" = "26.17"
    "E3" = "4.97%"
    "D4" = "5.087"
    "E4" = "0.80%"
    "D5" = "0.05604"
    "E5" = "-0.29%"
    "D6" = "6.482"
    "E6" = "-0.84%"
    "D7" = "0.8133"
    "E7" = "0.45%"
    "D8" = "0.8450"
    "E8" = "0.82%"
    "D9" = "0.02847"
    "E9" = "0.93%"
    "D10" = "0.09389"
    "E10" = "-0.24%"
    "D11" = "0.001519"
    "E11" = "0.69%"
    "D12" = "0.006118"
    "E12" = "-2.42%"
    "D13" = "3.600"
    "E13" = "2.88%"
    "D14" = "3.013"
    "E14" = "0.78%"
    "D15" = "2.055"
    "E15" = "0.08%"
    "D16" = "0.009899"
    "E16" = "1,553.41%"
    "E17" = "0.70%"
    "D18" = "0.1339"
    "E18" = "0.24%"
    "D19" = "0.07008"
    "D20" = "0.03204"
    "E20" = "-2.88%"
    "E21" = "0.47%"
    "D22" = "3.739"
    "E22" = "0.10%"
    "D23" = "0.04661"
    "E23" = "-0.25%"
    "E24" = "-1.38%"
    "D25" = "0.001247"
    "E25" = "0.40%"
    "D26" = "0.004584"
    "E26" = "1.26%"
    "D27" = "0.00009603"
    "D28" = "0.0001938"
    "E28" = "-0.03%"
    "E40" = "1.13%"
    "D41" = "0.006172"
    "E41" = "-1.04%"
    "D42" = "0.1056"
    "E42" = "0.41%"
    "D43" = "0.002501"
    "E43" = "-7.33%"
    "D44" = "0.008774"
    "E44" = "4.93%"
    "D45" = "0.00005297"
    "E45" = "0.65%"
    "E46" = "0.08%"
    "E47" = "-39.95%"
    "D48" = "0.002659"
    "E48" = "29.81%"
    "E49" = "0.08%"
    "E50" = "0.08%"
}

foreach ($key in $textCells.Keys) {
    $ws.Range($key).Value = $textCells[$key]
}

foreach ($key in $numericTextCells.Keys) {
    $r = $ws.Range($key)
    $r.NumberFormat = "@"
    $r.Value = $numericTextCells[$key]
}

Write-Output "Applied $($textCells.Count) text cell updates and $($numericTextCells.Count) numeric-text cell updates."
